$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Replace the old validation-message strings (now obsolete) with the new
# failure-forward JSP path on sheet "SalvarCliente".
$ws1.Range("H3").Value = "cliente/falha-cadastrar-cliente.jsp"
$ws1.Range("H4").Value = "cliente/falha-cadastrar-cliente.jsp"
$ws1.Range("H5").Value = "cliente/falha-cadastrar-cliente.jsp"
$ws1.Range("H6").Value = "cliente/falha-cadastrar-cliente.jsp"
$ws1.Range("H7").Value = "cliente/falha-cadastrar-cliente.jsp"
$ws1.Range("H8").Value = "cliente/falha-cadastrar-cliente.jsp"

# Underline the passing "Pobre" result to mark the happy-path outcome.
$ws1.Range("H9").Font.Underline = $true

# Move the active selection on sheet1 to H2 (where it landed after testing).
$ws1.Range("H2").Select()
